$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,7).Value = 0.131499
$ws.Cells.Item(2,8).Value = 0.394497
$ws.Cells.Item(2,9).Value = 0.3654391092296077
$ws.Cells.Item(2,10).Value = 0.3654391092296077
$ws.Cells.Item(2,13).Value = 14.440165
$ws.Cells.Item(2,14).Value = 43.320495
$ws.Cells.Item(2,15).Value = 0.1441015470002482
$ws.Cells.Item(2,16).Value = 0.1441015470002482
$ws.Cells.Item(2,17).Value = 1.898867257335
$ws.Cells.Item(2,18).Value = 17.089805316015
$ws.Cells.Item(2,19).Value = 0.05266034097437917
$ws.Cells.Item(2,20).Value = 0.05266034097437915

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,7).Value = 0.131499
$ws.Cells.Item(3,8).Value = 0.394497
$ws.Cells.Item(3,9).Value = 0.3654391092296077
$ws.Cells.Item(3,10).Value = 0.3654391092296077
$ws.Cells.Item(3,13).Value = 38.54369466666667
$ws.Cells.Item(3,14).Value = 115.631084
$ws.Cells.Item(3,15).Value = 0.3846359116098663
$ws.Cells.Item(3,16).Value = 0.3846359116098662
$ws.Cells.Item(3,17).Value = 5.068457304972
$ws.Cells.Item(3,18).Value = 45.616115744748
$ws.Cells.Item(3,19).Value = 0.1405610049164276
$ws.Cells.Item(3,20).Value = 0.1405610049164276

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,7).Value = 0.131499
$ws.Cells.Item(4,8).Value = 0.394497
$ws.Cells.Item(4,9).Value = 0.3654391092296077
$ws.Cells.Item(4,10).Value = 0.3654391092296077
$ws.Cells.Item(4,13).Value = 21.954262
$ws.Cells.Item(4,14).Value = 65.862786
$ws.Cells.Item(4,15).Value = 0.2190863551385157
$ws.Cells.Item(4,16).Value = 0.2190863551385156
$ws.Cells.Item(4,17).Value = 2.886963498738
$ws.Cells.Item(4,18).Value = 25.982671488642
$ws.Cells.Item(4,19).Value = 0.08006272246618065
$ws.Cells.Item(4,20).Value = 0.08006272246618064

# Row 5
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,4).Value = "Resolving-Mac"
$ws.Cells.Item(5,7).Value = 0.131499
$ws.Cells.Item(5,8).Value = 0.394497
$ws.Cells.Item(5,9).Value = 0.3654391092296077
$ws.Cells.Item(5,10).Value = 0.3654391092296077
$ws.Cells.Item(5,13).Value = 25.27013633333333
$ws.Cells.Item(5,14).Value = 75.81040899999999
$ws.Cells.Item(5,15).Value = 0.2521761862513699
$ws.Cells.Item(5,16).Value = 0.2521761862513699
$ws.Cells.Item(5,17).Value = 3.322997657697
$ws.Cells.Item(5,18).Value = 29.906978919273
$ws.Cells.Item(5,19).Value = 0.09215504087262028
$ws.Cells.Item(5,20).Value = 0.09215504087262028

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,7).Value = 0.2283393333333333
$ws.Cells.Item(6,8).Value = 0.685018
$ws.Cells.Item(6,9).Value = 0.6345608907703922
$ws.Cells.Item(6,10).Value = 0.6345608907703922
$ws.Cells.Item(6,13).Value = 14.440165
$ws.Cells.Item(6,14).Value = 43.320495
$ws.Cells.Item(6,15).Value = 0.1441015470002482
$ws.Cells.Item(6,16).Value = 0.1441015470002482
$ws.Cells.Item(6,17).Value = 3.297257649323333
$ws.Cells.Item(6,18).Value = 29.67531884391
$ws.Cells.Item(6,19).Value = 0.09144120602586907
$ws.Cells.Item(6,20).Value = 0.09144120602586905

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,7).Value = 0.2283393333333333
$ws.Cells.Item(7,8).Value = 0.685018
$ws.Cells.Item(7,9).Value = 0.6345608907703922
$ws.Cells.Item(7,10).Value = 0.6345608907703922
$ws.Cells.Item(7,13).Value = 38.54369466666667
$ws.Cells.Item(7,14).Value = 115.631084
$ws.Cells.Item(7,15).Value = 0.3846359116098663
$ws.Cells.Item(7,16).Value = 0.3846359116098662
$ws.Cells.Item(7,17).Value = 8.801041544390223
$ws.Cells.Item(7,18).Value = 79.209373899512
$ws.Cells.Item(7,19).Value = 0.2440749066934386
$ws.Cells.Item(7,20).Value = 0.2440749066934386

# Row 8
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,4).Value = "MuSCs"
$ws.Cells.Item(8,7).Value = 0.2283393333333333
$ws.Cells.Item(8,8).Value = 0.685018
$ws.Cells.Item(8,9).Value = 0.6345608907703922
$ws.Cells.Item(8,10).Value = 0.6345608907703922
$ws.Cells.Item(8,13).Value = 21.954262
$ws.Cells.Item(8,14).Value = 65.862786
$ws.Cells.Item(8,15).Value = 0.2190863551385157
$ws.Cells.Item(8,16).Value = 0.2190863551385156
$ws.Cells.Item(8,17).Value = 5.013021548905334
$ws.Cells.Item(8,18).Value = 45.117193940148
$ws.Cells.Item(8,19).Value = 0.139023632672335
$ws.Cells.Item(8,20).Value = 0.139023632672335

# Row 9
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,4).Value = "Resolving-Mac"
$ws.Cells.Item(9,7).Value = 0.2283393333333333
$ws.Cells.Item(9,8).Value = 0.685018
$ws.Cells.Item(9,9).Value = 0.6345608907703922
$ws.Cells.Item(9,10).Value = 0.6345608907703922
$ws.Cells.Item(9,13).Value = 25.27013633333333
$ws.Cells.Item(9,14).Value = 75.81040899999999
$ws.Cells.Item(9,15).Value = 0.2521761862513699
$ws.Cells.Item(9,16).Value = 0.2521761862513699
$ws.Cells.Item(9,17).Value = 5.770166083595777
$ws.Cells.Item(9,18).Value = 51.931494752362
$ws.Cells.Item(9,19).Value = 0.1600211453787496
$ws.Cells.Item(9,20).Value = 0.1600211453787496
